$wb = $excel.ActiveWorkbook

# ---- Sheet "Means" ----
$ws1 = $wb.Worksheets.Item("Means")

# New header cells for the 5-mile and 10-mile radius columns
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New column values (F = within 5 miles, G = within 10 miles) for each variable row
$ws1.Range("F2").Value = 71
$ws1.Range("G2").Value = 62

$ws1.Range("F3").Value = 27
$ws1.Range("G3").Value = 35

$ws1.Range("F4").Value = 2.4
$ws1.Range("G4").Value = 3.4

$ws1.Range("F5").Value = 5.7
$ws1.Range("G5").Value = 4.5

$ws1.Range("F6").Value = 51
$ws1.Range("G6").Value = 45

$ws1.Range("F7").Value = 12
$ws1.Range("G7").Value = 13

$ws1.Range("F8").Value = 8.8
$ws1.Range("G8").Value = 7.7

# Row 9 (Total Cancer Risk) existing B:E values updated, plus new F/G
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 34
$ws1.Range("D9").Value = 50
$ws1.Range("E9").Value = 50
$ws1.Range("F9").Value = 52
$ws1.Range("G9").Value = 49

# Row 10 (Total Respiratory) existing B:E values updated, plus new F/G
$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.47
$ws1.Range("D10").Value = 0.6
$ws1.Range("E10").Value = 0.6
$ws1.Range("F10").Value = 0.58
$ws1.Range("G10").Value = 0.54

# ---- Sheet "Standard Deviations" ----
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

$ws2.Range("F2").Value = 33
$ws2.Range("G2").Value = 31

$ws2.Range("F3").Value = 35
$ws2.Range("G3").Value = 33

$ws2.Range("F4").Value = 3.8
$ws2.Range("G4").Value = 4.4

$ws2.Range("F5").Value = 9.2
$ws2.Range("G5").Value = 7.6

$ws2.Range("F6").Value = 28
$ws2.Range("G6").Value = 21

$ws2.Range("F7").Value = 11
$ws2.Range("G7").Value = 14

$ws2.Range("F8").Value = 13
$ws2.Range("G8").Value = 9.7

# Row 9 (Total Cancer Risk SD) existing B:E values updated, plus new F/G
$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 5.5
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 4.7
$ws2.Range("G9").Value = 4.9

# Row 10 (Total Respiratory SD) existing C value updated, plus new F/G
$ws2.Range("C10").Value = 0.056
$ws2.Range("F10").Value = 0.083
$ws2.Range("G10").Value = 0.064
